$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 9465
$ws.Range("D2").Value = 8384
$ws.Range("E2").Value = 0.8857897517168516
$ws.Range("F2").Value = 0.8845748048111416
$ws.Range("G2").Value = 0.09679007923257567
$ws.Range("H2").Value = 0.08561806544481057
$ws.Range("I2").Value = 41250664.96178105
$ws.Range("J2").Value = 14430500.42754652
$ws.Range("L2").Value = 14430500.42754652
$ws.Range("M2").Value = 55681165.38932757
$ws.Range("N2").Value = 801737041.7372
$ws.Range("O2").Value = 784037234.7332001
$ws.Range("P2").Value = 0.0179990441707403
$ws.Range("Q2").Value = 0.018405376413605
$ws.Range("C3").Value = 9648
$ws.Range("D3").Value = 8554
$ws.Range("E3").Value = 0.886608623548922
$ws.Range("F3").Value = 0.8849575832816057
$ws.Range("G3").Value = 0.101949910499804
$ws.Range("H3").Value = 0.09022134641168254
$ws.Range("I3").Value = 48011620.77939813
$ws.Range("J3").Value = 17541061.32943826
$ws.Range("L3").Value = 17541061.32943826
$ws.Range("M3").Value = 65552682.10883638
$ws.Range("N3").Value = 837547624.534428
$ws.Range("O3").Value = 820067448.510358
$ws.Range("P3").Value = 0.02094335989453603
$ws.Range("Q3").Value = 0.02138977880575723
$ws.Range("C4").Value = 9842
$ws.Range("D4").Value = 8741
$ws.Range("E4").Value = 0.8881324933956513
$ws.Range("F4").Value = 0.8866910123757354
$ws.Range("G4").Value = 0.1061872908635677
$ws.Range("H4").Value = 0.09415531643725347
$ws.Range("I4").Value = 54607314.52462393
$ws.Range("J4").Value = 20518986.07138735
$ws.Range("L4").Value = 20518986.07138735
$ws.Range("M4").Value = 75126300.5960113
$ws.Range("N4").Value = 875630888.769419
$ws.Range("O4").Value = 858181940.7634652
$ws.Range("P4").Value = 0.02343337396448407
$ws.Range("Q4").Value = 0.02390983204928902
$ws.Range("C5").Value = 10033
$ws.Range("D5").Value = 8884
$ws.Range("E5").Value = 0.8854779228545799
$ws.Range("F5").Value = 0.8836284066043366
$ws.Range("G5").Value = 0.1096290891937461
$ws.Range("H5").Value = 0.09687137740175458
$ws.Range("I5").Value = 60627454.10565276
$ws.Range("J5").Value = 23195982.15371279
$ws.Range("L5").Value = 23195982.15371279
$ws.Range("M5").Value = 83823436.25936554
$ws.Range("N5").Value = 914110715.330657
$ws.Range("O5").Value = 896624609.8676838
$ws.Range("P5").Value = 0.02537546247373571
$ws.Range("Q5").Value = 0.02587033848773776
$ws.Range("C6").Value = 10231
$ws.Range("D6").Value = 9096
$ws.Range("E6").Value = 0.889062652722119
$ws.Range("F6").Value = 0.8870684610883558
$ws.Range("G6").Value = 0.1086865141011976
$ws.Range("H6").Value = 0.09641237880480726
$ws.Range("I6").Value = 63928479.04605511
$ws.Range("J6").Value = 24475449.96408697
$ws.Range("L6").Value = 24475449.96408697
$ws.Range("M6").Value = 88403929.01014209
$ws.Range("N6").Value = 955558181.6979581
$ws.Range("O6").Value = 937966355.8247766
$ws.Range("P6").Value = 0.0256137725916342
$ws.Range("Q6").Value = 0.02609416618420722
